$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B23").Value = 6307
$ws.Range("C23").Value = 998
$ws.Range("D23").Value = 5874319
$ws.Range("E23").Value = 931.3967020770573
$ws.Range("F23").Value = 8.218943033630755
$ws.Range("G23").Value = 3.850156087408951
$ws.Range("H23").Value = 25.87162006234101
